# SouthPlatte-WaterEvents.xlsx regression-test data update.
#
# Adds an "IconOriginalName" / "IconName" pair of lookup columns to the
# EventTypes sheet (columns C & D), mapping each event type to the icon
# file that represents it, and makes EventTypes the active sheet/tab
# (it was EventData before).

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("EventTypes")

# --- New column headers (row 1) --------------------------------------
$ws2.Range("C1").Value = "IconOriginalName"

# --- Column C: original icon file names --------------------------------
# Written in this order so the shared-string table comes out exactly as
# it did in the source workbook (the column was originally populated
# from an external, differently-ordered icon catalogue).
$ws2.Range("C13").Value = "cloud_rain_icon.png"
$ws2.Range("C6").Value  = "water_sea.png"
$ws2.Range("C5").Value  = "sun.png"
$ws2.Range("C4").Value  = "document_text_edit_32.png"
$ws2.Range("C3").Value  = "document_text_add_32.png"
$ws2.Range("C10").Value = "blur_gear.png"
$ws2.Range("C8").Value  = "snow_flake.png"
$ws2.Range("C11").Value = "group_half_32.png"
$ws2.Range("C12").Value = "oilwell.png"
$ws2.Range("C9").Value  = "document_text_information_32.png"
$ws2.Range("C2").Value  = "Basic-Scales-of-Balance-icon.png"
$ws2.Range("C7").Value  = "chat_exclamation.png"

# --- Column D: normalized "<EventType>.png" icon names ------------------
$ws2.Range("D2").Value  = "ColoradoLawCourtCase.png"
$ws2.Range("D1").Value  = "IconName"
$ws2.Range("D3").Value  = "ColoradoLawImplemented.png"
$ws2.Range("D4").Value  = "ColoradoLawPassed.png"
$ws2.Range("D5").Value  = "Drought.png"
$ws2.Range("D6").Value  = "Flood.png"
$ws2.Range("D7").Value  = "HighGroundWaterReports.png"
$ws2.Range("D8").Value  = "Snowpack.png"
$ws2.Range("D9").Value  = "Study.png"
$ws2.Range("D10").Value = "Technology.png"
$ws2.Range("D11").Value = "WellAugmentationGroup.png"
$ws2.Range("D12").Value = "WellConstruction.png"
$ws2.Range("D13").Value = "WetYears.png"

# --- Column widths for the two new columns -----------------------------
$ws2.Columns.Item(3).ColumnWidth = 35.1
$ws2.Columns.Item(4).ColumnWidth = 31.6

# --- Make EventTypes the active sheet/tab (was EventData) --------------
$ws2.Activate()
$ws2.Range("D14").Select() | Out-Null
